# Remove the quest system entry from the MainIcon table.
# Row 7 (Id=6, Name="传记" / biography quest entry, Icon="MainIcon3") is deleted
# entirely. Excel will shift rows 8:21 up to 7:20, automatically updating the
# worksheet dimension, the attached table ("表1") range/autofilter, and
# compacting the shared-strings table to drop now-unused strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("7:7").Delete()

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("C8").Select()
